$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-17 down to 8-18.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new weekly price record.
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44533
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112022
$ws.Range("G7").Value = "Arveja Verde"
$ws.Range("H7").Value = "Perfection"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14375
$ws.Range("N7").Value = "`$/malla 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 575
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
